# Change accuracy measure to informedness (balanced_accuracy)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hyperparameter text for Decision Tree (D2): clf__max_depth 7 -> 9
$ws.Range("D2").Value = "{'clf__criterion': 'gini', 'clf__max_depth': 9, 'clf__min_samples_leaf': 3, 'clf__min_samples_split': 0.05}"

# Update hyperparameter text for Random Forest (D5): criterion gini->entropy, min_samples_leaf 3->1
$ws.Range("D5").Value = "{'rf__criterion': 'entropy', 'rf__max_depth': 9, 'rf__min_samples_leaf': 1, 'rf__min_samples_split': 0.05}"

# Update Top Score column (C) with new balanced-accuracy-based scores
$ws.Range("C2").Value = 0.6980567603353105
$ws.Range("C3").Value = 0.80642373817349
$ws.Range("C4").Value = 0.9060831631487044
$ws.Range("C5").Value = 0.7733390130674143
$ws.Range("C6").Value = 0.7024704789650172
$ws.Range("C7").Value = 0.8580617694972213
